$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows to append (mirrors the style-less formatting already used by
# rows 16-20, i.e. no explicit cell style).
$newRows = @(
    @("gamodemy1@gmail.com", "20-June, 22:55", 4, 2, 0, 0),
    @("gamodemy1@gmail.com", "20-June, 23:04", 4, 2, 0, 0),
    @("gamodemy1@gmail.com", "20-June, 23:07", 4, 2, 0, 0)
)

$startRow = 21
$styleSourceRow = 16

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]

    for ($c = 1; $c -le 6; $c++) {
        # Pull formatting from a cell known to carry no explicit style so the
        # appended cells stay plain (matching the rest of the un-styled rows).
        $ws.Cells.Item($styleSourceRow, 1).Copy()
        $ws.Cells.Item($r, $c).PasteSpecial(-4122)
    }

    $ws.Cells.Item($r, 1).Value = $rowData[0]
    $ws.Cells.Item($r, 2).Value = $rowData[1]
    $ws.Cells.Item($r, 3).Value = $rowData[2]
    $ws.Cells.Item($r, 4).Value = $rowData[3]
    $ws.Cells.Item($r, 5).Value = $rowData[4]
    $ws.Cells.Item($r, 6).Value = $rowData[5]
}

$excel.CutCopyMode = 0
